# Apply the "Add files via upload" update to the COVID-19 vaccine priority
# populations tracker: refresh the as-of date, the national summary roll-up
# counts, and a handful of state rows (California, Mississippi, New Jersey)
# whose prioritization text was revised; also fix a stale footnote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title row: "Timeframe: as of ..." ---
$ws.Range("A2").Value = "Timeframe: as of January 14, 2021"

# --- United States roll-up row (row 4) ---
$ws.Range("C4").Value = "Follows ACIP Recommendations:12; Varies from ACIP Recommendations: 32; Not updated: 7"
$ws.Range("F4").Value = "Phase 1a: 31; Phase 1b: 15; Phase determined by counties: 5"

# --- California (row 9) ---
$ws.Range("F9").Value = "Phase 1b statewide: administering to people ages 65+ only. Counties determine when to include additional populations"

# --- Mississippi (row 29) ---
$ws.Range("C29").Value = "Includes ACIP recommended groups plus: people ages 65+; people ages 16-64 with high-risk medical conditions"
$ws.Range("D29").Value = "Includes all essential workers"
$ws.Range("F29").Value = "Phase 1b; administering to people ages 65+; people ages 16-64 with high-risk medical conditions only"

# --- New Jersey (row 35) ---
$ws.Range("C35").Value = "Includes ACIP recommended groups plus; people ages 65+; people ages 16-64 with high-risk medical conditions"
$ws.Range("D35").Value = "Includes other essential workers"
$ws.Range("F35").Value = "Phase 1b; administering to people ages 65+; people ages 16-64 with high-risk medical conditions; and fire/law enforcement personnel only"

# --- Footnotes: row 77 ("3. Massachusetts...") now mirrors the Tennessee footnote text ---
$ws.Range("A77").Value = "5. Tennessee has proposed two additional phases; phase 2a/b includes people ages 55-64 and critical infrastructure workers; phase 3 includes people ages 45-54, people living in congregate settings, grocery workers, and residents of correctional facilities"

# --- Reflect the cell that was active/selected when the file was last saved ---
$ws.Range("I9").Select()
